# Streamline the wellplate import template so the "Sample" column is a
# generic, empty column (instead of a pre-filled "sample_ID" column with
# placeholder IDs) — this makes the template reusable for research plan
# table imports as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column B header from "sample_ID" to "Sample"
$ws.Range("B1").Value = "Sample"

# Remove the placeholder sample IDs (10001, 10002, ...) from column B,
# leaving the cells blank while keeping their existing formatting.
$ws.Range("B2:B97").ClearContents()

# Reflect the author's last selection in the saved file.
$ws.Range("C11").Select() | Out-Null
